# Rename ID-style headers to lowercase, per commit:
# "changed IDs to lowercase id, ammended and tidied all structural checks,
#  some work on pkgdown vignettes"

$wb = $excel.ActiveWorkbook

# --- studies sheet: study_ID -> study_id ---
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("A1").Value = "study_id"

# --- surveys sheet: survey_ID -> survey_id, lat -> latitude, lon -> longitude ---
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"

# --- update the active selection on each sheet to match the new view state ---
$wsStudies.Range("C13").Select() | Out-Null
$wsSurveys.Range("F2").Select() | Out-Null

$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("B5").Select() | Out-Null

# prevalence becomes the final active sheet/tab
$wsPrevalence = $wb.Worksheets.Item("prevalence")
$wsPrevalence.Range("B6").Select() | Out-Null
